$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 103
$ws.Range("A20:E20").Copy()
$ws.Range("A103:E103").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H103:Q103").PasteSpecial(-4122)
$ws.Range("A103").Value = 43339.45383068287
$ws.Range("B103").Value = "Quadrat survey"
$ws.Range("C103").Value = "Little Birch Lake"
$ws.Range("D103").Value = 43334.0
$ws.Range("E103").Value = "Aislyn, Austen"
$ws.Range("H103").Value = "0.5m x 0.5m"
$ws.Range("I103").Value = 18.0
$ws.Range("J103").Value = 11.0
$ws.Range("K103").Value = 45.76223
$ws.Range("L103").Value = -94.79833
$ws.Range("M103").Value = 275.0
$ws.Range("N103").Value = 0.0014699074090458453
$ws.Range("O103").Value = 0.0012384259243845008
$ws.Range("P103").Value = 0.0116666666654055
$ws.Range("Q103").Value = 2.0

# Row 104
$ws.Range("A20:E20").Copy()
$ws.Range("A104:E104").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H104:Q104").PasteSpecial(-4122)
$ws.Range("A104").Value = 43339.45461554398
$ws.Range("B104").Value = "Quadrat survey"
$ws.Range("C104").Value = "Little Birch Lake"
$ws.Range("D104").Value = 43334.0
$ws.Range("E104").Value = "Aislyn, Austen"
$ws.Range("H104").Value = "0.5m x 0.5m"
$ws.Range("I104").Value = 22.0
$ws.Range("J104").Value = 10.0
$ws.Range("K104").Value = 45.76801
$ws.Range("L104").Value = -94.80041
$ws.Range("M104").Value = 300.0
$ws.Range("N104").Value = 0.0015046296321088448
$ws.Range("O104").Value = 0.0011921296318178065
$ws.Range("P104").Value = 0.010833333333721384
$ws.Range("Q104").Value = 2.0

# Row 105
$ws.Range("A20:E20").Copy()
$ws.Range("A105:E105").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H105:Q105").PasteSpecial(-4122)
$ws.Range("A105").Value = 43339.4560662963
$ws.Range("B105").Value = "Quadrat survey"
$ws.Range("C105").Value = "Little Birch Lake"
$ws.Range("D105").Value = 43334.0
$ws.Range("E105").Value = "Aislyn, Austen"
$ws.Range("H105").Value = "0.5m x 0.5m"
$ws.Range("I105").Value = 12.0
$ws.Range("J105").Value = 9.0
$ws.Range("K105").Value = 45.77292
$ws.Range("L105").Value = -94.80136
$ws.Range("M105").Value = 315.0
$ws.Range("N105").Value = 0.001111111108912155
$ws.Range("O105").Value = 0.0008333333316841163
$ws.Range("P105").Value = 0.008784722223936114
$ws.Range("Q105").Value = 2.0

# Row 106
$ws.Range("A20:E20").Copy()
$ws.Range("A106:E106").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H106:Q106").PasteSpecial(-4122)
$ws.Range("A106").Value = 43339.45690800926
$ws.Range("B106").Value = "Quadrat survey"
$ws.Range("C106").Value = "Little Birch Lake"
$ws.Range("D106").Value = 43334.0
$ws.Range("E106").Value = "Aislyn, Austen"
$ws.Range("H106").Value = "0.5m x 0.5m"
$ws.Range("I106").Value = 10.0
$ws.Range("J106").Value = 8.0
$ws.Range("K106").Value = 45.77744
$ws.Range("L106").Value = -94.79878
$ws.Range("M106").Value = 310.0
$ws.Range("N106").Value = 0.0009259259240934625
$ws.Range("O106").Value = 0.000787037039117422
$ws.Range("P106").Value = 0.008217592592700385
$ws.Range("Q106").Value = 2.0

# Row 107
$ws.Range("A20:E20").Copy()
$ws.Range("A107:E107").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H107:Q107").PasteSpecial(-4122)
$ws.Range("A107").Value = 43339.45793721065
$ws.Range("B107").Value = "Quadrat survey"
$ws.Range("C107").Value = "Little Birch Lake"
$ws.Range("D107").Value = 43334.0
$ws.Range("E107").Value = "Aislyn, Austen"
$ws.Range("H107").Value = "0.5m x 0.5m"
$ws.Range("I107").Value = 14.0
$ws.Range("J107").Value = 7.0
$ws.Range("K107").Value = 45.78306
$ws.Range("L107").Value = -94.79498
$ws.Range("M107").Value = 285.0
$ws.Range("N107").Value = 0.001331018516793847
$ws.Range("O107").Value = 0.0008680555547471158
$ws.Range("P107").Value = 0.010277777779265307
$ws.Range("Q107").Value = 2.0

# Row 108
$ws.Range("A20:E20").Copy()
$ws.Range("A108:E108").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H108:Q108").PasteSpecial(-4122)
$ws.Range("A108").Value = 43339.45941582176
$ws.Range("B108").Value = "Quadrat survey"
$ws.Range("C108").Value = "Little Birch Lake"
$ws.Range("D108").Value = 43334.0
$ws.Range("E108").Value = "Aislyn, Austen"
$ws.Range("H108").Value = "0.5m x 0.5m"
$ws.Range("I108").Value = 20.0
$ws.Range("J108").Value = 3.0
$ws.Range("K108").Value = 45.78365
$ws.Range("L108").Value = -94.78806
$ws.Range("M108").Value = 60.0
$ws.Range("N108").Value = 0.0014351851859828457
$ws.Range("O108").Value = 0.0009722222239361145
$ws.Range("P108").Value = 0.013981481482915115
$ws.Range("Q108").Value = 2.0

# Row 109
$ws.Range("A20:E20").Copy()
$ws.Range("A109:E109").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H109:Q109").PasteSpecial(-4122)
$ws.Range("A109").Value = 43339.460224131944
$ws.Range("B109").Value = "Quadrat survey"
$ws.Range("C109").Value = "Little Birch Lake"
$ws.Range("D109").Value = 43334.0
$ws.Range("E109").Value = "Aislyn, Austen"
$ws.Range("H109").Value = "0.5m x 0.5m"
$ws.Range("I109").Value = 26.0
$ws.Range("J109").Value = 4.0
$ws.Range("K109").Value = 45.77947
$ws.Range("L109").Value = -94.78256
$ws.Range("M109").Value = 320.0
$ws.Range("N109").Value = 0.001678240740147885
$ws.Range("O109").Value = 0.0012962962937308475
$ws.Range("P109").Value = 0.012766203704813961
$ws.Range("Q109").Value = 2.0

# Row 110
$ws.Range("A20:E20").Copy()
$ws.Range("A110:E110").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H110:Q110").PasteSpecial(-4122)
$ws.Range("A110").Value = 43339.46091677083
$ws.Range("B110").Value = "Quadrat survey"
$ws.Range("C110").Value = "Little Birch Lake"
$ws.Range("D110").Value = 43334.0
$ws.Range("E110").Value = "Aislyn, Austen"
$ws.Range("H110").Value = "0.5m x 0.5m"
$ws.Range("I110").Value = 20.0
$ws.Range("J110").Value = 5.0
$ws.Range("K110").Value = 45.78445
$ws.Range("L110").Value = -94.77834
$ws.Range("M110").Value = 250.0
$ws.Range("N110").Value = 0.0011689814855344594
$ws.Range("O110").Value = 0.000891203701030463
$ws.Range("P110").Value = 0.010625000002619345
$ws.Range("Q110").Value = 2.0

# Row 111
$ws.Range("A20:E20").Copy()
$ws.Range("A111:E111").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H111:Q111").PasteSpecial(-4122)
$ws.Range("A111").Value = 43339.4617322338
$ws.Range("B111").Value = "Quadrat survey"
$ws.Range("C111").Value = "Little Birch Lake"
$ws.Range("D111").Value = 43334.0
$ws.Range("E111").Value = "Aislyn, Austen"
$ws.Range("H111").Value = "0.5m x 0.5m"
$ws.Range("I111").Value = 28.0
$ws.Range("J111").Value = 6.0
$ws.Range("K111").Value = 45.79258
$ws.Range("L111").Value = -94.78363
$ws.Range("M111").Value = 245.0
$ws.Range("N111").Value = 0.0015509259246755391
$ws.Range("O111").Value = 0.0011805555550381541
$ws.Range("P111").Value = 0.014004629629198462
$ws.Range("Q111").Value = 2.0

# Row 112
$ws.Range("A20:E20").Copy()
$ws.Range("A112:E112").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H112:Q112").PasteSpecial(-4122)
$ws.Range("A112").Value = 43339.46246686342
$ws.Range("B112").Value = "Quadrat survey"
$ws.Range("C112").Value = "Little Birch Lake"
$ws.Range("D112").Value = 43334.0
$ws.Range("E112").Value = "Aislyn, Austen"
$ws.Range("H112").Value = "0.5m x 0.5m"
$ws.Range("I112").Value = 14.0
$ws.Range("J112").Value = 1.0
$ws.Range("K112").Value = 45.79403
$ws.Range("L112").Value = -94.79369
$ws.Range("M112").Value = 120.0
$ws.Range("N112").Value = 0.0013541666703531519
$ws.Range("O112").Value = 0.0011342592624714598
$ws.Range("P112").Value = 0.006377314814017154
$ws.Range("Q112").Value = 0.75

# Row 113
$ws.Range("A20:E20").Copy()
$ws.Range("A113:E113").PasteSpecial(-4122)
$ws.Range("H20:Q20").Copy()
$ws.Range("H113:Q113").PasteSpecial(-4122)
$ws.Range("A113").Value = 43339.46452400463
$ws.Range("B113").Value = "Quadrat survey"
$ws.Range("C113").Value = "Little Birch Lake"
$ws.Range("D113").Value = 43334.0
$ws.Range("E113").Value = "Aislyn, Austen"
$ws.Range("H113").Value = "0.5m x 0.5m"
$ws.Range("I113").Value = 20.0
$ws.Range("J113").Value = 2.0
$ws.Range("K113").Value = 45.78797
$ws.Range("L113").Value = -94.79337
$ws.Range("M113").Value = 100.0
$ws.Range("N113").Value = 0.001157407408754807
$ws.Range("O113").Value = 0.0010648148163454607
$ws.Range("P113").Value = 0.011481481480586808
$ws.Range("Q113").Value = 2.0

# Row 114
$ws.Range("A102:E102").Copy()
$ws.Range("A114:E114").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F114:G114").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I114:Q114").PasteSpecial(-4122)
$ws.Range("A114").Value = 43339.46731096065
$ws.Range("B114").Value = "Double observer no distance"
$ws.Range("C114").Value = "Little Birch Lake"
$ws.Range("D114").Value = 43333.0
$ws.Range("E114").Value = "Aislyn, Austen"
$ws.Range("F114").Value = "Austen"
$ws.Range("G114").Value = "Aislyn"
$ws.Range("I114").Value = 13.0
$ws.Range("J114").Value = 12.0
$ws.Range("K114").Value = 45.76497
$ws.Range("L114").Value = -94.80565
$ws.Range("M114").Value = 60.0
$ws.Range("N114").Value = 0.004027777773444541
$ws.Range("O114").Value = 0.0011689814855344594
$ws.Range("P114").Value = 0.007581018522614613
$ws.Range("Q114").Value = 1.5

# Row 115
$ws.Range("A102:E102").Copy()
$ws.Range("A115:E115").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F115:G115").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I115:Q115").PasteSpecial(-4122)
$ws.Range("A115").Value = 43339.46915204861
$ws.Range("B115").Value = "Double observer no distance"
$ws.Range("C115").Value = "Little Birch Lake"
$ws.Range("D115").Value = 43333.0
$ws.Range("E115").Value = "Aislyn, Austen"
$ws.Range("F115").Value = "Aislyn"
$ws.Range("G115").Value = "Austen"
$ws.Range("I115").Value = 19.0
$ws.Range("J115").Value = 11.0
$ws.Range("K115").Value = 45.76175
$ws.Range("L115").Value = -94.79804
$ws.Range("M115").Value = 270.0
$ws.Range("N115").Value = 0.004583333335176576
$ws.Range("O115").Value = 0.0017939814788405783
$ws.Range("P115").Value = 0.02071759258979
$ws.Range("Q115").Value = 1.5

# Row 116
$ws.Range("A102:E102").Copy()
$ws.Range("A116:E116").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F116:G116").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I116:Q116").PasteSpecial(-4122)
$ws.Range("A116").Value = 43339.47008313658
$ws.Range("B116").Value = "Double observer no distance"
$ws.Range("C116").Value = "Little Birch Lake"
$ws.Range("D116").Value = 43333.0
$ws.Range("E116").Value = "Aislyn, Austen"
$ws.Range("F116").Value = "Austen"
$ws.Range("G116").Value = "Aislyn"
$ws.Range("I116").Value = 20.0
$ws.Range("J116").Value = 10.0
$ws.Range("K116").Value = 45.7687
$ws.Range("L116").Value = -94.80112
$ws.Range("M116").Value = 290.0
$ws.Range("N116").Value = 0.0036111111112404615
$ws.Range("O116").Value = 0.00148148147854954
$ws.Range("P116").Value = 0.0071064814837882295
$ws.Range("Q116").Value = 2.0

# Row 117
$ws.Range("A102:E102").Copy()
$ws.Range("A117:E117").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F117:G117").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I117:Q117").PasteSpecial(-4122)
$ws.Range("A117").Value = 43339.47086568287
$ws.Range("B117").Value = "Double observer no distance"
$ws.Range("C117").Value = "Little Birch Lake"
$ws.Range("D117").Value = 43333.0
$ws.Range("E117").Value = "Aislyn, Austen"
$ws.Range("F117").Value = "Aislyn"
$ws.Range("G117").Value = "Austen"
$ws.Range("I117").Value = 12.0
$ws.Range("J117").Value = 9.0
$ws.Range("K117").Value = 45.77284
$ws.Range("L117").Value = -94.80145
$ws.Range("M117").Value = 330.0
$ws.Range("N117").Value = 0.003298611110949423
$ws.Range("O117").Value = 0.0010532407395658083
$ws.Range("P117").Value = 0.01608796296204673
$ws.Range("Q117").Value = 2.0

# Row 118
$ws.Range("A102:E102").Copy()
$ws.Range("A118:E118").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F118:G118").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I118:Q118").PasteSpecial(-4122)
$ws.Range("A118").Value = 43339.471509375
$ws.Range("B118").Value = "Double observer no distance"
$ws.Range("C118").Value = "Little Birch Lake"
$ws.Range("D118").Value = 43333.0
$ws.Range("E118").Value = "Aislyn, Austen"
$ws.Range("F118").Value = "Austen"
$ws.Range("G118").Value = "Aislyn"
$ws.Range("I118").Value = 11.0
$ws.Range("J118").Value = 8.0
$ws.Range("K118").Value = 45.77736
$ws.Range("L118").Value = -94.7988
$ws.Range("M118").Value = 340.0
$ws.Range("N118").Value = 0.0030902777798473835
$ws.Range("O118").Value = 0.000787037039117422
$ws.Range("P118").Value = 0.016018518515920732
$ws.Range("Q118").Value = 2.0

# Row 119
$ws.Range("A102:E102").Copy()
$ws.Range("A119:E119").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F119:G119").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I119:Q119").PasteSpecial(-4122)
$ws.Range("A119").Value = 43339.472419259255
$ws.Range("B119").Value = "Double observer no distance"
$ws.Range("C119").Value = "Little Birch Lake"
$ws.Range("D119").Value = 43333.0
$ws.Range("E119").Value = "Aislyn, Austen"
$ws.Range("F119").Value = "Aislyn"
$ws.Range("G119").Value = "Austen"
$ws.Range("I119").Value = 16.0
$ws.Range("J119").Value = 7.0
$ws.Range("K119").Value = 45.78196
$ws.Range("L119").Value = -94.79467
$ws.Range("M119").Value = 275.0
$ws.Range("N119").Value = 0.003078703703067731
$ws.Range("O119").Value = 0.0014467592554865405
$ws.Range("P119").Value = 0.0196180555576575
$ws.Range("Q119").Value = 2.0

# Row 120
$ws.Range("A102:E102").Copy()
$ws.Range("A120:E120").PasteSpecial(-4122)
$ws.Range("F102:G102").Copy()
$ws.Range("F120:G120").PasteSpecial(-4122)
$ws.Range("I102:Q102").Copy()
$ws.Range("I120:Q120").PasteSpecial(-4122)
$ws.Range("A120").Value = 43339.473133657404
$ws.Range("B120").Value = "Double observer no distance"
$ws.Range("C120").Value = "Little Birch Lake"
$ws.Range("D120").Value = 43333.0
$ws.Range("E120").Value = "Aislyn, Austen"
$ws.Range("F120").Value = "Austen"
$ws.Range("G120").Value = "Aislyn"
$ws.Range("I120").Value = 24.0
$ws.Range("J120").Value = 3.0
$ws.Range("K120").Value = 45.78364
$ws.Range("L120").Value = -94.78804
$ws.Range("M120").Value = 90.0
$ws.Range("N120").Value = 0.004293981481168885
$ws.Range("O120").Value = 0.0015046296321088448
$ws.Range("P120").Value = 0.036585648151231
$ws.Range("Q120").Value = 2.0

$excel.CutCopyMode = $false
